$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update F5 formatting: mark it red (style change s="1" -> s="5") ---
$ws.Range("F5").Font.Color = 255

# --- Append two new test rows (33 and 34) below the existing table data ---

# Row 33: copy formatting from row 32 (same column layout), then set values
$ws.Range("A32:F32").Copy()
$ws.Range("A33:F33").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("A33").Value = 32
$ws.Range("B33").Value = "Comportement des option d'analyse"
$ws.Range("C33").Value = "Changement des chiffres significatifs"
$ws.Range("D33").Value = "./Data/*.csv"
$ws.Range("E33").Value = "Passer l'option des chiffres significatifs sur 5. Configurer et analyser le sample. Exporter les graphiques individuels. Ouvrir la fenêtre des variables."
$ws.Range("F33").Value = "Les chiffres significatifs devraient être de 5 partout où cela est possible: Les tables des graphiques exportés; Fenêtre des variables; Fenêtre post-analyse et tableaux exportés (VarWin & PAWin)."
$ws.Rows(33).RowHeight = 57.6

# Row 34: same formatting/group ("Comportement des option d'analyse"), different task
$ws.Range("A32:F32").Copy()
$ws.Range("A34:F34").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("A34").Value = 33
$ws.Range("B34").Value = "Comportement des option d'analyse"
$ws.Range("C34").Value = "Changement dedéfinition limite élastique"
$ws.Range("D34").Value = "./Data/*.csv"
$ws.Range("E34").Value = "Passer l'option de Définition de limite élastique à 0.5% pour optenir Rp0.5 par la suite. Configurer et analyser le sample. Exporter les graphiques individuels."
$ws.Range("F34").Value = "La définition de Re passe à 0.5%. Les résultats calculés pour Re doivent avoir été adaptés en conséquence. Vérifier sur les graphique ""Contrainte - Déformation"" que la droite change bien de position."
$ws.Rows(34).RowHeight = 57.6

# --- Update sheet view: clear the scrolled top-left cell and move selection to G6 ---
[void]$ws.Range("G6").Select()
